$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text format
# first, otherwise Excel auto-converts the typed string into a numeric value and the
# original text (including fixed decimal places) would not be preserved.
$textCells = @("D5", "D6", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D42", "D44", "D46", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the diff
$ws.Range("D2").Value = "64.879.56"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.154.83"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "571.10"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("D6").Value = "150.72"
$ws.Range("E6").Value = "  +5.90%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.151.45"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("E9").Value = "  +5.12%  "
$ws.Range("E10").Value = "  +5.76%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "0.505"
$ws.Range("E12").Value = "  +7.78%  "
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").Value = "  +15.03%  "
$ws.Range("D14").Value = "38.21"
$ws.Range("E14").Value = "  +10.20%  "
$ws.Range("D15").Value = "3.672.52"
$ws.Range("D16").Value = "64.976.43"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").Value = "  +7.93%  "
$ws.Range("D18").Value = "3.157.64"
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "516.25"
$ws.Range("E20").Value = "  +8.34%  "
$ws.Range("D21").Value = "14.93"
$ws.Range("E21").Value = "  +7.14%  "
$ws.Range("D22").Value = "0.737"
$ws.Range("E22").Value = "  +9.67%  "
$ws.Range("D23").Value = "15.43"
$ws.Range("E23").Value = "  +10.07%  "
$ws.Range("E24").Value = "  +5.22%  "
$ws.Range("D25").Value = "85.07"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.93"
$ws.Range("E27").Value = "  +5.70%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "8.94"
$ws.Range("E28").Value = "  +13.02%  "
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +8.00%  "
$ws.Range("D30").Value = "27.87"
$ws.Range("E30").Value = "  +6.94%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "2.72"
$ws.Range("E32").Value = "  +11.02%  "
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +4.55%  "
$ws.Range("D34").Value = "6.22"
$ws.Range("E34").Value = "  +11.65%  "
$ws.Range("D35").Value = "6.65"
$ws.Range("E35").Value = "  +8.39%  "
$ws.Range("D36").Value = "55.90"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "485.79"
$ws.Range("E37").Value = "  +10.53%  "
$ws.Range("E38").Value = "  +8.44%  "
$ws.Range("E39").Value = "  +5.07%  "
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").Value = "  +5.14%  "
$ws.Range("D41").Value = "3.119.49"
$ws.Range("E41").Value = "  +5.82%  "
$ws.Range("D42").Value = "8.68"
$ws.Range("E42").Value = "  +6.49%  "
$ws.Range("E43").Value = "  +5.34%  "
$ws.Range("D44").Value = "0.293"
$ws.Range("E44").Value = "  +14.05%  "
$ws.Range("E45").Value = "  +18.27%  "
$ws.Range("D46").Value = "29.74"
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("D47").Value = "0.0₃0579"
$ws.Range("E47").Value = "  +13.12%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  +12.84%  "
$ws.Range("D51").Value = "121.10"
$ws.Range("E51").Value = "  +3.90%  "
